# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in column G ("Recorded By") of the active worksheet, for every row where
# that exact value is present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Host "Rows updated:" $changed
